$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "b,1,1,1"
$ws.Range("A3").Value = "b"
$ws.Range("B3").Value = ""
